# Generate Report for Handback
# Update the "Correspond Handback DateTime" timestamps (column G) on each
# language sheet to reflect the new report-generation run. Only cells that
# currently hold the sheet's previous "last generated" timestamp are
# touched; any other (older / exception) timestamps in that column are
# left untouched.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "zh-cn"; Old = "2016-02-22 08:48:32"; New = "2016-02-22 08:49:28" },
    @{ Sheet = "de-de"; Old = "2016-02-22 08:48:45"; New = "2016-02-22 08:49:38" },
    @{ Sheet = "ja-jp"; Old = "2016-02-22 08:48:57"; New = "2016-02-22 08:49:48" },
    @{ Sheet = "zh-tw"; Old = "2016-02-22 08:49:07"; New = "2016-02-22 08:49:58" }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 7)
        if ($cell.Value2 -eq $u.Old) {
            $cell.Value = $u.New
        }
    }
}
